$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.009.20'
$ws.Range("E2").Value = '  +1.40%  '
$ws.Range("D3").Value = '2.397.96'
$ws.Range("E3").Value = '  +0.91%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''509.16'
$ws.Range("E5").Value = '  +1.96%  '
$ws.Range("D6").Value = '''134.99'
$ws.Range("E6").Value = '  +4.69%  '
$ws.Range("D7").Value = '''0.998'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '''0.555'
$ws.Range("E8").Value = '  +0.40%  '
$ws.Range("D9").Value = '2.404.97'
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("D10").Value = '''0.0986'
$ws.Range("E10").Value = '  +3.29%  '
$ws.Range("D11").Value = '''0.151'
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("D12").Value = '''0.337'
$ws.Range("E12").Value = '  +6.55%  '
$ws.Range("D13").Value = '''4.69'
$ws.Range("E13").Value = '  +1.07%  '
$ws.Range("D14").Value = '2.818.56'
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("D15").Value = '56.927.33'
$ws.Range("E15").Value = '  +1.41%  '
$ws.Range("D16").Value = '''21.97'
$ws.Range("E16").Value = '  +2.43%  '
$ws.Range("D17").Value = '''0.0000134'
$ws.Range("E17").Value = '  +2.42%  '
$ws.Range("D18").Value = '2.373.32'
$ws.Range("E18").Value = '  -1.43%  '
$ws.Range("D19").Value = '''10.20'
$ws.Range("E19").Value = '  +1.24%  '
$ws.Range("D20").Value = '''4.07'
$ws.Range("E20").Value = '  +1.37%  '
$ws.Range("D21").Value = '''311.85'
$ws.Range("E21").Value = '  +0.91%  '
$ws.Range("D22").Value = '''6.24'
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("D24").Value = '''5.64'
$ws.Range("E24").Value = '  +1.03%  '
$ws.Range("D25").Value = '''65.47'
$ws.Range("E25").Value = '  +1.62%  '
$ws.Range("D26").Value = '''0.998'
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").Value = '''0.375'
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("D28").Value = '''0.152'
$ws.Range("E28").Value = '  +1.96%  '
$ws.Range("D29").Value = '''7.40'
$ws.Range("E29").Value = '  +2.04%  '
$ws.Range("D30").Value = '''173.00'
$ws.Range("E30").Value = '  +0.21%  '
$ws.Range("D31").Value = '0.0₃0734'
$ws.Range("E31").Value = '  +3.20%  '
$ws.Range("D32").Value = '''1.66'
$ws.Range("E32").Value = '  +0.61%  '
$ws.Range("D33").Value = '''1.14'
$ws.Range("E33").Value = '  +4.07%  '
$ws.Range("D34").Value = '''5.89'
$ws.Range("E34").Value = '  -3.38%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").Value = '''0.996'
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").Value = '''17.93'
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("D38").Value = '''1.21'
$ws.Range("E38").Value = '  +0.23%  '
$ws.Range("D39").Value = '''3.89'
$ws.Range("E39").Value = '  +3.28%  '
$ws.Range("D40").Value = '''36.80'
$ws.Range("E40").Value = '  +2.73%  '
$ws.Range("D41").Value = '''0.821'
$ws.Range("E41").Value = '  +4.86%  '
$ws.Range("D42").Value = '''1.45'
$ws.Range("E42").Value = '  +1.82%  '
$ws.Range("D43").Value = '''132.17'
$ws.Range("E43").Value = '  +2.52%  '
$ws.Range("D44").Value = '''3.41'
$ws.Range("E44").Value = '  +2.48%  '
$ws.Range("D45").Value = '''4.94'
$ws.Range("E45").Value = '  +3.50%  '
$ws.Range("D46").Value = '''0.570'
$ws.Range("E46").Value = '  +1.76%  '
$ws.Range("D47").Value = '''0.0913'
$ws.Range("E47").Value = '  +1.57%  '
$ws.Range("D48").Value = '''250.57'
$ws.Range("E48").Value = '  -0.10%  '
$ws.Range("D49").Value = '''0.0489'
$ws.Range("E49").Value = '  +1.10%  '
$ws.Range("D50").Value = '''0.0211'
$ws.Range("E50").Value = '  +2.16%  '
$ws.Range("D51").Value = '''17.36'
$ws.Range("E51").Value = '  +7.35%  '
